# "add exam to database and retrieve exam"
# The exams table gains a new "subject" column (inserted between exam_name
# and teacher_id) and the exam_name / subject values are refreshed to match
# what is now stored in the database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("subject"); everything from the old C onward
# (teacher_id, total_marks, start_time, end_time) shifts one column right.
$ws.Columns("C").Insert()

# Header row
$ws.Range("C1").Value = "subject"

# Refreshed exam_name values (col B) pulled from the database
$ws.Range("B2").Value = "UT01"
$ws.Range("B3").Value = "UT01"
$ws.Range("B4").Value = "ut02"
$ws.Range("B5").Value = "ut01"
$ws.Range("B6").Value = "ut03"

# New subject values (col C) pulled from the database
$ws.Range("C2").Value = "english"
$ws.Range("C3").Value = "physics"
$ws.Range("C4").Value = "match"
$ws.Range("C5").Value = "chemisrty"
$ws.Range("C6").Value = "ip"

# The old exam_name strings forced word-wrap onto two lines (extra row
# height); the new short codes fit on one line, so let the rows shrink
# back to the default height.
$ws.Rows("2:6").AutoFit() | Out-Null

# Column sizing to fit the refreshed table layout.
$ws.Columns("B").ColumnWidth = 10.333333333333334
$ws.Columns("C").ColumnWidth = 10.5
$ws.Columns("D").ColumnWidth = 10.333333333333334
$ws.Columns("E").ColumnWidth = 11.666666666666666

# Restore the active selection used when the sheet was last saved.
$ws.Range("J11").Select() | Out-Null
